# Add two new columns (I: I0, J: IF) to the worksheet, mirroring the
# existing header/style pattern used by columns B..H.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header formatting already used by the other header cells
# (bold font, thin box border, centered / top-aligned).
$headerRange = $ws.Range("I1:J1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4160     # xlTop
$headerRange.Borders.LineStyle = 1         # xlContinuous
$headerRange.Borders.Weight = 2            # xlThin

# --- Data rows (2-10) ---
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 2

$ws.Range("I3").Value = 8
$ws.Range("J3").Value = 8

$ws.Range("I4").Value = 8
$ws.Range("J4").Value = 8

$ws.Range("I5").Value = 8
$ws.Range("J5").Value = 8

$ws.Range("I6").Value = 8
$ws.Range("J6").Value = 8

$ws.Range("I7").Value = 9
$ws.Range("J7").Value = 9

$ws.Range("I8").Value = 9
$ws.Range("J8").Value = 9

$ws.Range("I9").Value = 9
$ws.Range("J9").Value = 9

$ws.Range("I10").Value = 5
$ws.Range("J10").Value = 5
